# Insert a new weekly price record as row 139 (pushing the existing rows
# 139-201 down to 140-202), matching the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 139..201 down by one to make room for the new record.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new data point.
$ws.Cells.Item(139, 1).Value = 3
$ws.Cells.Item(139, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(139, 3).Value = "Coquimbo"
$ws.Cells.Item(139, 4).Value = 44460
$ws.Cells.Item(139, 5).Value = 5
$ws.Cells.Item(139, 6).Value = 100112040
$ws.Cells.Item(139, 7).Value = "Cilantro"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 160
$ws.Cells.Item(139, 11).Value = 2500
$ws.Cells.Item(139, 12).Value = 2500
$ws.Cells.Item(139, 13).Value = 2500
$ws.Cells.Item(139, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(139, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(139, 16).Value = 833
$ws.Cells.Item(139, 17).Value = 3
$ws.Cells.Item(139, 18).Value = "Hortaliza"
